$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1854838709677419
$ws.Range("C2").Value = 0.5833333333333334
$ws.Range("J2").Value = 0.005376344086021506
$ws.Range("P2").Value = 0.1290322580645161
$ws.Range("S2").Value = 0.09677419354838709
$ws.Range("B3").Value = 0.008771929824561403
$ws.Range("C3").Value = 0.01754385964912281
$ws.Range("J3").Value = 0.03070175438596491
$ws.Range("P3").Value = 0.6973684210526315
$ws.Range("S3").Value = 0.2456140350877193
$ws.Range("J4").Value = 0.06557377049180328
$ws.Range("P4").Value = 0.6065573770491803
$ws.Range("S4").Value = 0.3278688524590164
$ws.Range("B6").Value = 0.08438818565400844
$ws.Range("D6").Value = 0.01265822784810127
$ws.Range("F6").Value = 0.06329113924050633
$ws.Range("J6").Value = 0.2278481012658228
$ws.Range("O6").Value = 0.02109704641350211
$ws.Range("Q6").Value = 0.1392405063291139
$ws.Range("R6").Value = 0.04219409282700422
$ws.Range("S6").Value = 0.4092827004219409
$ws.Range("B7").Value = 0.1260162601626016
$ws.Range("D7").Value = 0.04065040650406504
$ws.Range("F7").Value = 0.04065040650406504
$ws.Range("J7").Value = 0.1463414634146341
$ws.Range("O7").Value = 0.01219512195121951
$ws.Range("Q7").Value = 0.1544715447154472
$ws.Range("R7").Value = 0.06504065040650407
$ws.Range("S7").Value = 0.4146341463414634
$ws.Range("B8").Value = 0.115546218487395
$ws.Range("D8").Value = 0.01470588235294118
$ws.Range("F8").Value = 0.06092436974789916
$ws.Range("J8").Value = 0.1218487394957983
$ws.Range("O8").Value = 0.02100840336134454
$ws.Range("Q8").Value = 0.1428571428571428
$ws.Range("R8").Value = 0.06722689075630252
$ws.Range("S8").Value = 0.4558823529411765
$ws.Range("B9").Value = 0.08383233532934131
$ws.Range("D9").Value = 0.04191616766467066
$ws.Range("F9").Value = 0.0658682634730539
$ws.Range("J9").Value = 0.1017964071856287
$ws.Range("O9").Value = 0.03592814371257485
$ws.Range("Q9").Value = 0.155688622754491
$ws.Range("R9").Value = 0.05988023952095808
$ws.Range("S9").Value = 0.4550898203592814
$ws.Range("B10").Value = 0.1248196248196248
$ws.Range("D10").Value = 0.0266955266955267
$ws.Range("F10").Value = 0.07431457431457432
$ws.Range("J10").Value = 0.09884559884559885
$ws.Range("O10").Value = 0.0202020202020202
$ws.Range("Q10").Value = 0.2012987012987013
$ws.Range("R10").Value = 0.07142857142857142
$ws.Range("S10").Value = 0.3823953823953824
$ws.Range("G11").Value = 0.1602067183462532
$ws.Range("J11").Value = 0.09560723514211886
$ws.Range("K11").Value = 0.2170542635658915
$ws.Range("L11").Value = 0.4987080103359173
$ws.Range("S11").Value = 0.02842377260981912
$ws.Range("G12").Value = 0.7222222222222222
$ws.Range("J12").Value = 0.2323232323232323
$ws.Range("K12").Value = 0.005050505050505051
$ws.Range("L12").Value = 0.0101010101010101
$ws.Range("S12").Value = 0.0303030303030303
$ws.Range("F13").Value = 0.01282051282051282
$ws.Range("G13").Value = 0.6025641025641025
$ws.Range("J13").Value = 0.2948717948717949
$ws.Range("S13").Value = 0.08974358974358974
$ws.Range("F15").Value = 0.02777777777777778
$ws.Range("H15").Value = 0.1296296296296296
$ws.Range("I15").Value = 0.06944444444444445
$ws.Range("J15").Value = 0.3564814814814815
$ws.Range("K15").Value = 0.05555555555555555
$ws.Range("M15").Value = 0.009259259259259259
$ws.Range("O15").Value = 0.04166666666666666
$ws.Range("S15").Value = 0.3101851851851852
$ws.Range("F16").Value = 0.02991452991452992
$ws.Range("H16").Value = 0.1752136752136752
$ws.Range("I16").Value = 0.05555555555555555
$ws.Range("J16").Value = 0.4017094017094017
$ws.Range("K16").Value = 0.1282051282051282
$ws.Range("M16").Value = 0.03846153846153846
$ws.Range("O16").Value = 0.03846153846153846
$ws.Range("S16").Value = 0.1324786324786325
$ws.Range("F17").Value = 0.01573033707865169
$ws.Range("H17").Value = 0.1573033707865168
$ws.Range("I17").Value = 0.0898876404494382
$ws.Range("J17").Value = 0.4292134831460674
$ws.Range("K17").Value = 0.1191011235955056
$ws.Range("M17").Value = 0.03370786516853932
$ws.Range("O17").Value = 0.04269662921348315
$ws.Range("S17").Value = 0.1123595505617977
$ws.Range("F18").Value = 0.03012048192771084
$ws.Range("H18").Value = 0.1927710843373494
$ws.Range("I18").Value = 0.07228915662650602
$ws.Range("J18").Value = 0.4036144578313253
$ws.Range("K18").Value = 0.1204819277108434
$ws.Range("M18").Value = 0.01807228915662651
$ws.Range("O18").Value = 0.06626506024096386
$ws.Range("S18").Value = 0.0963855421686747
$ws.Range("F19").Value = 0.01028101439342015
$ws.Range("H19").Value = 0.2090472926662097
$ws.Range("I19").Value = 0.05825908156271419
$ws.Range("J19").Value = 0.3858807402330363
$ws.Range("K19").Value = 0.1274845784784099
$ws.Range("M19").Value = 0.03495544893762851
$ws.Range("O19").Value = 0.06031528444139822
$ws.Range("S19").Value = 0.113776559287183

Write-Host "Applied all changes"
